$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.800.82"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "'1.879.68"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'324.42"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").Value = "'0.4677"
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("D8").Value = "'0.3936"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("D9").Value = "'0.07924"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "'0.9817"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").Value = "'22.38"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").Value = "'1.933.76"
$ws.Range("E12").Value = "  +6.29%  "
$ws.Range("D13").Value = "'5.742"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").Value = "'7.018"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "'0.06982"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "'88.73"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "'1.007"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "'16.97"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").Value = "'28.827.56"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").Value = "'5.349"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").Value = "'2.128"
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("D25").Value = "'2.118.21"
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("D26").Value = "'153.48"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").Value = "'19.40"
$ws.Range("E27").Value = "  +0.92%  "
$ws.Range("D28").Value = "'5.780"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "'119.89"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").Value = "'0.09402"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").Value = "'0.9370"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'5.312"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'1.357"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "'3.349"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "'0.05911"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("D37").Value = "'0.02130"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "'1.162"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").Value = "'7.894"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("D41").Value = "'0.1797"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "'10.01"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "'0.07314"
$ws.Range("E43").Value = "  +3.64%  "
$ws.Range("E44").Value = "  +2.02%  "
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("D46").Value = "'0.5358"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("D47").Value = "'1.847"
$ws.Range("D48").Value = "'113.98"
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("D49").Value = "'2.081"
$ws.Range("E49").Value = "  -6.39%  "
$ws.Range("D50").Value = "'2.377"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "'1.006"
$ws.Range("E51").Value = "  +0.55%  "
